# The commit rotates the data rows 2, 3 and 4 by one position:
#   new row 2 <- old row 3
#   new row 3 <- old row 4
#   new row 4 <- old row 2
# (row 1 is the header row and rows 5/6 are untouched by the commit).
#
# Only the columns that actually differ between the rows are touched
# (A, Q, R, S, Y, AA, AW, AX, AY), and only when the incoming value is
# actually different from what is already there, so every other cell -
# including the blank placeholder cells - is left completely undisturbed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "Q", "R", "S", "Y", "AA", "AW", "AX", "AY")

# Columns Y/AA hold dates stored as plain text (e.g. "2019-10-26"); force
# text format before writing so the values round-trip as text instead of
# being re-interpreted as date serial numbers.
$ws.Range("Y2:Y4").NumberFormat = "@"
$ws.Range("AA2:AA4").NumberFormat = "@"

# Snapshot the current (pre-rotation) values of each relevant column.
$before = @{}
foreach ($col in $cols) {
    $before[$col] = @{
        2 = $ws.Range($col + "2").Value2
        3 = $ws.Range($col + "3").Value2
        4 = $ws.Range($col + "4").Value2
    }
}

# Write the rotated values back: row2 <- row3, row3 <- row4, row4 <- row2.
# Skip the write whenever the target already holds the incoming value so
# cells that are not really changing keep their original representation.
foreach ($col in $cols) {
    $new2 = $before[$col][3]
    $new3 = $before[$col][4]
    $new4 = $before[$col][2]

    if ($ws.Range($col + "2").Value2 -ne $new2) { $ws.Range($col + "2").Value2 = $new2 }
    if ($ws.Range($col + "3").Value2 -ne $new3) { $ws.Range($col + "3").Value2 = $new3 }
    if ($ws.Range($col + "4").Value2 -ne $new4) { $ws.Range($col + "4").Value2 = $new4 }
}
